$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 3276
$ws.Cells.Item(40, 9).Value = 3245
$ws.Cells.Item(40, 10).Value = 3400
$ws.Cells.Item(40, 11).Value = 3245
$ws.Cells.Item(40, 12).Value = 3400
$ws.Cells.Item(40, 13).Value = -3070
$ws.Cells.Item(40, 14).Value = -3750
$ws.Cells.Item(41, 8).Value = 6944725
$ws.Cells.Item(41, 9).Value = 8547312
$ws.Cells.Item(41, 10).Value = 180.66667
$ws.Cells.Item(41, 11).Value = 8547312
$ws.Cells.Item(41, 12).Value = 180.66667
$ws.Cells.Item(41, 13).Value = -8546872
$ws.Cells.Item(41, 14).Value = -1060.66667
$ws.Cells.Item(64, 8).Value = 3910649.5
$ws.Cells.Item(64, 9).Value = 8336778.5
$ws.Cells.Item(64, 10).Value = 5241.294
$ws.Cells.Item(64, 11).Value = 8336778.5
$ws.Cells.Item(64, 12).Value = 5241.294
$ws.Cells.Item(64, 13).Value = -8336530.5
$ws.Cells.Item(64, 14).Value = -5737.294
$ws.Cells.Item(67, 8).Value = 3910649.5
$ws.Cells.Item(67, 9).Value = 8336778.5
$ws.Cells.Item(67, 10).Value = 5241.294
$ws.Cells.Item(67, 11).Value = 8336778.5
$ws.Cells.Item(67, 12).Value = 5241.294
$ws.Cells.Item(67, 13).Value = -8335920.5
$ws.Cells.Item(67, 14).Value = -6957.294
$ws.Cells.Item(74, 8).Value = 3444
$ws.Cells.Item(74, 9).Value = 3298.6667
$ws.Cells.Item(74, 10).Value = 3589.3333
$ws.Cells.Item(74, 11).Value = 3298.6667
$ws.Cells.Item(74, 12).Value = 3589.3333
$ws.Cells.Item(74, 13).Value = -2362.6667
$ws.Cells.Item(74, 14).Value = -5461.3333
$ws.Cells.Item(77, 8).Value = 3444
$ws.Cells.Item(77, 9).Value = 3298.6667
$ws.Cells.Item(77, 10).Value = 3589.3333
$ws.Cells.Item(77, 11).Value = 16493.3335
$ws.Cells.Item(77, 12).Value = 17946.6665
$ws.Cells.Item(77, 13).Value = -11813.3335
$ws.Cells.Item(77, 14).Value = -27306.6665
$ws.Cells.Item(138, 8).Value = 6423828.5
$ws.Cells.Item(138, 9).Value = 2553867
$ws.Cells.Item(138, 10).Value = 8200204.5
$ws.Cells.Item(138, 11).Value = 7661601
$ws.Cells.Item(138, 12).Value = 24600613.5
$ws.Cells.Item(138, 13).Value = -7656461
$ws.Cells.Item(138, 14).Value = -24610893.5
$ws.Cells.Item(141, 8).Value = 2402.3333
$ws.Cells.Item(141, 9).Value = 2402.3333
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 11).Value = 7206.999899999999
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 14).Value = -2026.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 1209.7
$ws.Cells.Item(4, 9).Value = 1209.7
$ws.Cells.Item(4, 11).Value = 1209.7
$ws.Cells.Item(4, 13).Value = -1093.7
$ws.Cells.Item(5, 8).Value = 111582.89
$ws.Cells.Item(5, 9).Value = 143435.42
$ws.Cells.Item(5, 10).Value = 99
$ws.Cells.Item(5, 11).Value = 143435.42
$ws.Cells.Item(5, 12).Value = 99
$ws.Cells.Item(5, 13).Value = -143323.42
$ws.Cells.Item(5, 14).Value = -323
$ws.Cells.Item(32, 8).Value = 21414.44
$ws.Cells.Item(32, 9).Value = 5790.403
$ws.Cells.Item(32, 10).Value = 91201.8
$ws.Cells.Item(32, 11).Value = 5790.403
$ws.Cells.Item(32, 12).Value = 91201.8
$ws.Cells.Item(32, 13).Value = -5503.403
$ws.Cells.Item(32, 14).Value = -91775.8
$ws.Cells.Item(63, 8).Value = 12600.833
$ws.Cells.Item(63, 9).Value = 14001.25
$ws.Cells.Item(63, 11).Value = 14001.25
$ws.Cells.Item(63, 13).Value = -13315.25
$ws.Cells.Item(66, 8).Value = 12600.833
$ws.Cells.Item(66, 9).Value = 14001.25
$ws.Cells.Item(66, 11).Value = 70006.25
$ws.Cells.Item(66, 13).Value = -66574.25
$ws.Cells.Item(74, 8).Value = 4788.75
$ws.Cells.Item(74, 9).Value = 965.0454999999999
$ws.Cells.Item(74, 10).Value = 10797.429
$ws.Cells.Item(74, 11).Value = 965.0454999999999
$ws.Cells.Item(74, 12).Value = 10797.429
$ws.Cells.Item(74, 13).Value = -91.04549999999995
$ws.Cells.Item(74, 14).Value = -12545.429
$ws.Cells.Item(77, 8).Value = 4788.75
$ws.Cells.Item(77, 9).Value = 965.0454999999999
$ws.Cells.Item(77, 10).Value = 10797.429
$ws.Cells.Item(77, 11).Value = 4825.2275
$ws.Cells.Item(77, 12).Value = 53987.145
$ws.Cells.Item(77, 13).Value = -457.2275
$ws.Cells.Item(77, 14).Value = -62723.145
$ws.Cells.Item(132, 8).Value = 3610.5789
$ws.Cells.Item(132, 9).Value = 2652.9167
$ws.Cells.Item(132, 10).Value = 5252.2856
$ws.Cells.Item(132, 11).Value = 7958.750100000001
$ws.Cells.Item(132, 12).Value = 15756.8568
$ws.Cells.Item(132, 13).Value = -5428.750100000001
$ws.Cells.Item(132, 14).Value = -20816.8568
$ws.Cells.Item(135, 8).Value = 55800
$ws.Cells.Item(135, 10).Value = 55800
$ws.Cells.Item(135, 12).Value = 55800
$ws.Cells.Item(135, 14).Value = -65940

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 111582.89
$ws.Cells.Item(4, 9).Value = 143435.42
$ws.Cells.Item(4, 10).Value = 99
$ws.Cells.Item(4, 11).Value = 143435.42
$ws.Cells.Item(4, 12).Value = 99
$ws.Cells.Item(4, 13).Value = -143320.42
$ws.Cells.Item(4, 14).Value = -329
$ws.Cells.Item(15, 8).Value = 4000
$ws.Cells.Item(15, 10).Value = 4000
$ws.Cells.Item(15, 12).Value = 4000
$ws.Cells.Item(15, 14).Value = -4454
$ws.Cells.Item(19, 8).Value = 50000
$ws.Cells.Item(19, 9).Value = 50000
$ws.Cells.Item(19, 11).Value = 50000
$ws.Cells.Item(19, 13).Value = -49827
$ws.Cells.Item(82, 8).Value = 21064.777
$ws.Cells.Item(82, 9).Value = 8369.166999999999
$ws.Cells.Item(82, 10).Value = 46456
$ws.Cells.Item(82, 11).Value = 8369.166999999999
$ws.Cells.Item(82, 12).Value = 46456
$ws.Cells.Item(82, 13).Value = -7986.166999999999
$ws.Cells.Item(82, 14).Value = -47222
$ws.Cells.Item(85, 8).Value = 21064.777
$ws.Cells.Item(85, 9).Value = 8369.166999999999
$ws.Cells.Item(85, 10).Value = 46456
$ws.Cells.Item(85, 11).Value = 8369.166999999999
$ws.Cells.Item(85, 12).Value = 46456
$ws.Cells.Item(85, 13).Value = -7043.166999999999
$ws.Cells.Item(85, 14).Value = -49108

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 1000000
$ws.Cells.Item(4, 10).Value = 1000000
$ws.Cells.Item(4, 12).Value = 1000000
$ws.Cells.Item(4, 14).Value = -1000224
$ws.Cells.Item(7, 8).Value = 61.42857
$ws.Cells.Item(7, 9).Value = 61
$ws.Cells.Item(7, 10).Value = 62.5
$ws.Cells.Item(7, 11).Value = 61
$ws.Cells.Item(7, 12).Value = 62.5
$ws.Cells.Item(7, 13).Value = 52
$ws.Cells.Item(7, 14).Value = -288.5
$ws.Cells.Item(31, 8).Value = 6140.625
$ws.Cells.Item(31, 10).Value = 6140.625
$ws.Cells.Item(31, 12).Value = 6140.625
$ws.Cells.Item(31, 14).Value = -6730.625
$ws.Cells.Item(34, 8).Value = 6140.625
$ws.Cells.Item(34, 10).Value = 6140.625
$ws.Cells.Item(34, 12).Value = 6140.625
$ws.Cells.Item(34, 14).Value = -6544.625
$ws.Cells.Item(56, 8).Value = 12833.333
$ws.Cells.Item(56, 9).Value = 2000
$ws.Cells.Item(56, 10).Value = 15000
$ws.Cells.Item(56, 11).Value = 2000
$ws.Cells.Item(56, 12).Value = 15000
$ws.Cells.Item(56, 13).Value = -1155
$ws.Cells.Item(56, 14).Value = -16690
$ws.Cells.Item(58, 8).Value = 26318004
$ws.Cells.Item(58, 9).Value = 41668188
$ws.Cells.Item(58, 10).Value = 3401.7856
$ws.Cells.Item(58, 11).Value = 41668188
$ws.Cells.Item(58, 12).Value = 3401.7856
$ws.Cells.Item(58, 13).Value = -41667985
$ws.Cells.Item(58, 14).Value = -3807.7856
$ws.Cells.Item(99, 8).Value = 10755304
$ws.Cells.Item(99, 9).Value = 2411.4783
$ws.Cells.Item(99, 11).Value = 2411.4783
$ws.Cells.Item(99, 13).Value = -913.4783000000002
$ws.Cells.Item(107, 8).Value = 987.7143
$ws.Cells.Item(107, 9).Value = 727.75
$ws.Cells.Item(107, 10).Value = 1334.3334
$ws.Cells.Item(107, 11).Value = 727.75
$ws.Cells.Item(107, 12).Value = 1334.3334
$ws.Cells.Item(107, 13).Value = 1192.25
$ws.Cells.Item(107, 14).Value = -5174.3334
$ws.Cells.Item(122, 8).Value = 1195.7826
$ws.Cells.Item(122, 9).Value = 1168.3182
$ws.Cells.Item(122, 11).Value = 3504.9546
$ws.Cells.Item(122, 13).Value = -1054.9546
$ws.Cells.Item(126, 8).Value = 10755304
$ws.Cells.Item(126, 9).Value = 2411.4783
$ws.Cells.Item(126, 11).Value = 7234.4349
$ws.Cells.Item(126, 13).Value = -4764.4349
$ws.Cells.Item(132, 8).Value = 3625080.8
$ws.Cells.Item(132, 9).Value = 4903248
$ws.Cells.Item(132, 11).Value = 14709744
$ws.Cells.Item(132, 13).Value = -14707214
$ws.Cells.Item(134, 8).Value = 31917240
$ws.Cells.Item(134, 9).Value = 41667830
$ws.Cells.Item(134, 10).Value = 21742710
$ws.Cells.Item(134, 11).Value = 125003490
$ws.Cells.Item(134, 12).Value = 65228130
$ws.Cells.Item(134, 13).Value = -125000955
$ws.Cells.Item(134, 14).Value = -65233200
$ws.Cells.Item(136, 8).Value = 26318004
$ws.Cells.Item(136, 9).Value = 41668188
$ws.Cells.Item(136, 10).Value = 3401.7856
$ws.Cells.Item(136, 11).Value = 125004564
$ws.Cells.Item(136, 12).Value = 10205.3568
$ws.Cells.Item(136, 13).Value = -125002014
$ws.Cells.Item(136, 14).Value = -15305.3568

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(69, 8).Value = 6666.6665
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(72, 8).Value = 6666.6665
$ws.Cells.Item(72, 9).Value = 0
$ws.Cells.Item(72, 11).Value = 0
$ws.Cells.Item(118, 8).Value = 1404.8334
$ws.Cells.Item(118, 9).Value = 885.8
$ws.Cells.Item(118, 11).Value = 2657.4
$ws.Cells.Item(118, 13).Value = -1414.4
$ws.Cells.Item(131, 8).Value = 16669029
$ws.Cells.Item(131, 10).Value = 19610462
$ws.Cells.Item(131, 12).Value = 58831386
$ws.Cells.Item(131, 14).Value = -58841466

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 1999.5
$ws.Cells.Item(113, 9).Value = 1999
$ws.Cells.Item(113, 10).Value = 2000
$ws.Cells.Item(113, 11).Value = 1999
$ws.Cells.Item(113, 12).Value = 2000
$ws.Cells.Item(113, 13).Value = 171
$ws.Cells.Item(113, 14).Value = -6340
$ws.Cells.Item(122, 8).Value = 2124.6206
$ws.Cells.Item(122, 9).Value = 2196.3845
$ws.Cells.Item(122, 11).Value = 6589.1535
$ws.Cells.Item(122, 13).Value = -4139.1535
$ws.Cells.Item(132, 8).Value = 3675.2856
$ws.Cells.Item(132, 9).Value = 3501.4546
$ws.Cells.Item(132, 10).Value = 3866.5
$ws.Cells.Item(132, 11).Value = 10504.3638
$ws.Cells.Item(132, 12).Value = 11599.5
$ws.Cells.Item(132, 13).Value = -7974.363799999999
$ws.Cells.Item(132, 14).Value = -16659.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(11, 8).Value = 4000
$ws.Cells.Item(11, 10).Value = 4000
$ws.Cells.Item(11, 12).Value = 4000
$ws.Cells.Item(11, 14).Value = -4280
$ws.Cells.Item(16, 8).Value = 7692882.5
$ws.Cells.Item(16, 9).Value = 9091416
$ws.Cells.Item(16, 10).Value = 950
$ws.Cells.Item(16, 11).Value = 9091416
$ws.Cells.Item(16, 12).Value = 950
$ws.Cells.Item(16, 13).Value = -9091246
$ws.Cells.Item(16, 14).Value = -1290
$ws.Cells.Item(46, 8).Value = 646
$ws.Cells.Item(46, 9).Value = 574.4
$ws.Cells.Item(46, 10).Value = 825
$ws.Cells.Item(46, 11).Value = 574.4
$ws.Cells.Item(46, 12).Value = 825
$ws.Cells.Item(46, 13).Value = -386.4
$ws.Cells.Item(46, 14).Value = -1201
$ws.Cells.Item(132, 8).Value = 4667.6665
$ws.Cells.Item(132, 9).Value = 3743.8572
$ws.Cells.Item(132, 10).Value = 5129.5713
$ws.Cells.Item(132, 11).Value = 11231.5716
$ws.Cells.Item(132, 12).Value = 15388.7139
$ws.Cells.Item(132, 13).Value = -8701.571599999999
$ws.Cells.Item(132, 14).Value = -20448.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1788.1765
$ws.Cells.Item(122, 9).Value = 1613.2667
$ws.Cells.Item(122, 11).Value = 4839.800099999999
$ws.Cells.Item(122, 13).Value = -2389.800099999999
$ws.Cells.Item(136, 8).Value = 2337.7256
$ws.Cells.Item(136, 9).Value = 617.8461
$ws.Cells.Item(136, 10).Value = 7927.3335
$ws.Cells.Item(136, 11).Value = 1853.5383
$ws.Cells.Item(136, 12).Value = 23782.0005
$ws.Cells.Item(136, 13).Value = 696.4617000000001
$ws.Cells.Item(136, 14).Value = -28882.0005

# Remove cells that no longer exist in the target (trailing cells dropped from rows)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(141, 14).ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(69, 13).ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(72, 13).ClearContents()
